# dimensionless numbers van Morgan
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- add a new (second) worksheet "Blad2" after Blad1 ---
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Blad2"
$ws.Activate()

# --- try to resize the saved window (best effort; cosmetic) ---
$wb.Windows.Item(1).Width = 24240
$wb.Windows.Item(1).Height = 13020

# --- new experiment "Morgan" data rows (14-18) ---
$ws.Range("A14").Value = "Morgan"
$ws.Range("B14").Value = "PEO600K"

$ws.Range("D14").Value = 44
$ws.Range("D15").Value = 44
$ws.Range("D16").Value = 44
$ws.Range("D17").Value = 44
$ws.Range("D18").Value = 44

$ws.Range("G14").Value = 0.05
$ws.Range("G15").Value = 0.1
$ws.Range("G16").Value = 0.2
$ws.Range("G17").Value = 0.5
$ws.Range("G18").Value = 1

$ws.Range("I14").Value = 0.00026
$ws.Range("I15").Value = 0.00041
$ws.Range("I16").Value = 0.00062
$ws.Range("I17").Value = 0.00136
$ws.Range("I18").Value = 0.00258

# the first new value keeps its own number format (5 decimals); rest keep
# the pre-existing style used by the surrounding empty cells
$ws.Range("I14").NumberFormat = "0.00000"

# --- Reynolds (K) column formula, filled down as one shared group ---
$ws.Range("K2:K18").Formula = "=I2*D2/(`$R`$1/1000)"

# K2 and K14 were (re)typed by hand in the source edit and stay as
# their own, non-shared formula entries
$ws.Range("K2").Formula = "=I2*D2/(`$R`$1/1000)"
$ws.Range("K11:K13").ClearContents()
$ws.Range("K14").Formula = "=I14*D14/(`$R`$1/1000)"

# column K now needs a fitted width, like I, J and L already have
$ws.Columns.Item(11).AutoFit()

# --- move the saved selection ---
$ws.Range("G21").Select()
